$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 829.8929000000001
$ws.Range("I80").Value = 1372.9
$ws.Range("J80").Value = 528.2222
$ws.Range("K80").Value = 4118.700000000001
$ws.Range("L80").Value = 1584.6666
$ws.Range("M80").Value = -3120.700000000001
$ws.Range("N80").Value = -3580.6666
$ws.Range("H83").Value = 829.8929000000001
$ws.Range("I83").Value = 1372.9
$ws.Range("J83").Value = 528.2222
$ws.Range("K83").Value = 12356.1
$ws.Range("L83").Value = 4753.999800000001
$ws.Range("M83").Value = -7364.1
$ws.Range("N83").Value = -14737.9998
$ws.Range("H138").Value = 9806235
$ws.Range("I138").Value = 17545998
$ws.Range("J138").Value = 2533.3333
$ws.Range("K138").Value = 52637994
$ws.Range("L138").Value = 7599.999899999999
$ws.Range("M138").Value = -52632854
$ws.Range("N138").Value = -17879.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10929.884
$ws.Range("I32").Value = 11420.308
$ws.Range("J32").Value = 9429.764999999999
$ws.Range("K32").Value = 11420.308
$ws.Range("L32").Value = 9429.764999999999
$ws.Range("M32").Value = -11133.308
$ws.Range("N32").Value = -10003.765
$ws.Range("H102").Value = 1374.9375
$ws.Range("I102").Value = 1428.5
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1428.5
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 193.5
$ws.Range("N102").Value = -4244
$ws.Range("H132").Value = 4311986
$ws.Range("I132").Value = 6580604
$ws.Range("J132").Value = 1611.75
$ws.Range("K132").Value = 19741812
$ws.Range("L132").Value = 4835.25
$ws.Range("M132").Value = -19739282
$ws.Range("N132").Value = -9895.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 526.7179599999999
$ws.Range("I94").Value = 418.5
$ws.Range("K94").Value = 418.5
$ws.Range("M94").Value = 32.5
$ws.Range("H99").Value = 1235.3636
$ws.Range("I99").Value = 1148.625
$ws.Range("J99").Value = 1466.6666
$ws.Range("K99").Value = 1148.625
$ws.Range("L99").Value = 1466.6666
$ws.Range("M99").Value = 349.375
$ws.Range("N99").Value = -4462.6666
$ws.Range("H134").Value = 2077.1516
$ws.Range("I134").Value = 1135.1522
$ws.Range("K134").Value = 3405.4566
$ws.Range("M134").Value = -870.4566

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 36850
$ws.Range("J20").Value = 36850
$ws.Range("L20").Value = 36850
$ws.Range("N20").Value = -37322
$ws.Range("H30").Value = 36850
$ws.Range("J30").Value = 36850
$ws.Range("L30").Value = 36850
$ws.Range("N30").Value = -37032
$ws.Range("H31").Value = 5212020.5
$ws.Range("I31").Value = 3849.1956
$ws.Range("J31").Value = 18521792
$ws.Range("K31").Value = 3849.1956
$ws.Range("L31").Value = 18521792
$ws.Range("M31").Value = -3554.1956
$ws.Range("N31").Value = -18522382
$ws.Range("H34").Value = 5212020.5
$ws.Range("I34").Value = 3849.1956
$ws.Range("J34").Value = 18521792
$ws.Range("K34").Value = 3849.1956
$ws.Range("L34").Value = 18521792
$ws.Range("M34").Value = -3647.1956
$ws.Range("N34").Value = -18522196
$ws.Range("H51").Value = 9598.799999999999
$ws.Range("J51").Value = 9598.799999999999
$ws.Range("L51").Value = 9598.799999999999
$ws.Range("N51").Value = -11070.8
$ws.Range("H58").Value = 1704.3334
$ws.Range("I58").Value = 1121.9
$ws.Range("J58").Value = 2600.3845
$ws.Range("K58").Value = 1121.9
$ws.Range("L58").Value = 2600.3845
$ws.Range("M58").Value = -918.9000000000001
$ws.Range("N58").Value = -3006.3845
$ws.Range("H61").Value = 9598.799999999999
$ws.Range("J61").Value = 9598.799999999999
$ws.Range("L61").Value = 9598.799999999999
$ws.Range("N61").Value = -10294.8
$ws.Range("H62").Value = 2076.524
$ws.Range("I62").Value = 2120.3333
$ws.Range("J62").Value = 1967
$ws.Range("K62").Value = 2120.3333
$ws.Range("L62").Value = 1967
$ws.Range("M62").Value = -1496.3333
$ws.Range("N62").Value = -3215
$ws.Range("H65").Value = 2076.524
$ws.Range("I65").Value = 2120.3333
$ws.Range("J65").Value = 1967
$ws.Range("K65").Value = 10601.6665
$ws.Range("L65").Value = 9835
$ws.Range("M65").Value = -7481.666499999999
$ws.Range("N65").Value = -16075
$ws.Range("H94").Value = 4301.8096
$ws.Range("I94").Value = 3330
$ws.Range("J94").Value = 4787.7144
$ws.Range("K94").Value = 3330
$ws.Range("L94").Value = 4787.7144
$ws.Range("M94").Value = -2879
$ws.Range("N94").Value = -5689.7144
$ws.Range("H107").Value = 408
$ws.Range("I107").Value = 426.65
$ws.Range("J107").Value = 361.375
$ws.Range("K107").Value = 426.65
$ws.Range("L107").Value = 361.375
$ws.Range("M107").Value = 1493.35
$ws.Range("N107").Value = -4201.375
$ws.Range("H128").Value = 36850
$ws.Range("J128").Value = 36850
$ws.Range("L128").Value = 36850
$ws.Range("N128").Value = -46810
$ws.Range("H132").Value = 8622168
$ws.Range("I132").Value = 10001344
$ws.Range("J132").Value = 2315.875
$ws.Range("K132").Value = 30004032
$ws.Range("L132").Value = 6947.625
$ws.Range("M132").Value = -30001502
$ws.Range("N132").Value = -12007.625
$ws.Range("H134").Value = 271504.4
$ws.Range("I134").Value = 961.7183
$ws.Range("J134").Value = 1401418
$ws.Range("K134").Value = 2885.1549
$ws.Range("L134").Value = 4204254
$ws.Range("M134").Value = -350.1549
$ws.Range("N134").Value = -4209324
$ws.Range("H136").Value = 1704.3334
$ws.Range("I136").Value = 1121.9
$ws.Range("J136").Value = 2600.3845
$ws.Range("K136").Value = 3365.7
$ws.Range("L136").Value = 7801.1535
$ws.Range("M136").Value = -815.7000000000003
$ws.Range("N136").Value = -12901.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9807215
$ws.Range("I80").Value = 18521170
$ws.Range("J80").Value = 4016.9375
$ws.Range("K80").Value = 18521170
$ws.Range("L80").Value = 4016.9375
$ws.Range("M80").Value = -18520172
$ws.Range("N80").Value = -6012.9375
$ws.Range("H83").Value = 9807215
$ws.Range("I83").Value = 18521170
$ws.Range("J83").Value = 4016.9375
$ws.Range("K83").Value = 92605850
$ws.Range("L83").Value = 20084.6875
$ws.Range("M83").Value = -92600858
$ws.Range("N83").Value = -30068.6875
$ws.Range("H122").Value = 5131135
$ws.Range("I122").Value = 13335873
$ws.Range("J122").Value = 3173.75
$ws.Range("K122").Value = 40007619
$ws.Range("L122").Value = 9521.25
$ws.Range("M122").Value = -40005169
$ws.Range("N122").Value = -14421.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5319.851
$ws.Range("I7").Value = 5186.778
$ws.Range("J7").Value = 5499.5
$ws.Range("K7").Value = 5186.778
$ws.Range("L7").Value = 5499.5
$ws.Range("M7").Value = -5074.778
$ws.Range("N7").Value = -5723.5
$ws.Range("H68").Value = 2180
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 1770
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 1770
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -3268
$ws.Range("H71").Value = 2180
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 1770
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 8850
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -16338
$ws.Range("H126").Value = 5319.851
$ws.Range("I126").Value = 5186.778
$ws.Range("J126").Value = 5499.5
$ws.Range("K126").Value = 15560.334
$ws.Range("L126").Value = 16498.5
$ws.Range("M126").Value = -13090.334
$ws.Range("N126").Value = -21438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2182
$ws.Range("I126").Value = 1401.7826
$ws.Range("J126").Value = 5771
$ws.Range("K126").Value = 4205.3478
$ws.Range("L126").Value = 17313
$ws.Range("M126").Value = -1735.3478
$ws.Range("N126").Value = -22253
$ws.Range("H132").Value = 1410.125
$ws.Range("I132").Value = 1402.075
$ws.Range("K132").Value = 4206.225
$ws.Range("M132").Value = -1676.225
$ws.Range("H136").Value = 841.6731
$ws.Range("I136").Value = 754.5263
$ws.Range("J136").Value = 1078.2142
$ws.Range("K136").Value = 2263.5789
$ws.Range("L136").Value = 3234.6426
$ws.Range("M136").Value = 286.4211
$ws.Range("N136").Value = -8334.642599999999
